$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.8
$ws.Range("I2").Value = 2.55
$ws.Range("L2").Value = 3.2
$ws.Range("AB2").Value = 29
$ws.Range("AM2").Value = 23

# Row 3 updates
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 5.5

# Row 4 updates
$ws.Range("G4").Value = 2.6
$ws.Range("I4").Value = 2.8
$ws.Range("J4").Value = 3.2
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 1.85
$ws.Range("AM4").Value = 29
$ws.Range("AN4").Value = 23
